$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "execute" column (C) for each test case was changed from "no" to "yes",
# so that every test is now marked to execute.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = "yes"
}
